# Speed modification to Crawler proccess
# Adds the "12 - 10" and "13 - 10" crawler run sheets (PullAndBear, Mango,
# Zara, Stradivarius, Bershka, MercedesCampuzano) after the existing
# "11 - 10 ..." sheets, and appends the newly-observed change timestamps to
# the tail end of the "11 - 10 ..." sheets that were still open when those
# later runs picked up more changes.

$wb = $excel.ActiveWorkbook

function New-CrawlerSheet {
    param(
        [string]$Name,
        [string[]]$Times
    )

    $count = $wb.Worksheets.Count
    $after = $wb.Worksheets.Item($count)
    $ws = $wb.Worksheets.Add($null, $after)
    $ws.Name = $Name

    $ws.Range("A1").Value2 = "Hora"
    $ws.Range("B1").Value2 = "Cambió"
    $ws.Range("C1").Value2 = "Nuevos"
    $ws.Range("D1").Value2 = "Actualizados"

    $header = $ws.Range("A1:D1")
    $header.Font.Bold = $true
    $header.Borders.LineStyle = 1
    $header.HorizontalAlignment = -4108
    $header.VerticalAlignment = -4160

    $r = 2
    foreach ($t in $Times) {
        $ws.Cells.Item($r, 1).Value2 = $t
        $ws.Cells.Item($r, 2).Value2 = $false
        $ws.Cells.Item($r, 3).Value2 = 0
        $ws.Cells.Item($r, 4).Value2 = 0
        $r = $r + 1
    }

    return $ws
}

# --- new sheets: "12 - 10" run ------------------------------------------
$null = New-CrawlerSheet "12 - 10 PullAndBear" @("23:28", "23:31")
$null = New-CrawlerSheet "12 - 10 Mango" @("23:28", "23:31")
$null = New-CrawlerSheet "12 - 10 Zara" @("23:28", "23:31")
$null = New-CrawlerSheet "12 - 10 Stradivarius" @("23:28", "23:31")
$null = New-CrawlerSheet "12 - 10 Bershka" @("23:28", "23:31")
$null = New-CrawlerSheet "12 - 10 MercedesCampuzano" @()

# --- new sheets: "13 - 10" run ------------------------------------------
$null = New-CrawlerSheet "13 - 10 MercedesCampuzano" @()
$null = New-CrawlerSheet "13 - 10 PullAndBear" @("0:6")
$null = New-CrawlerSheet "13 - 10 Mango" @("0:6")
$null = New-CrawlerSheet "13 - 10 Zara" @("0:6")
$null = New-CrawlerSheet "13 - 10 Stradivarius" @("0:6")
$null = New-CrawlerSheet "13 - 10 Bershka" @("0:6")

# --- append late-arriving change rows onto the "11 - 10" sheets --------

function Add-Row {
    param($ws, [int]$Row, [string]$Time)
    $ws.Cells.Item($Row, 1).Value2 = $Time
    $ws.Cells.Item($Row, 2).Value2 = $false
    $ws.Cells.Item($Row, 3).Value2 = 0
    $ws.Cells.Item($Row, 4).Value2 = 0
}

$wsPullAndBear = $wb.Worksheets.Item("11 - 10 PullAndBear")
Add-Row $wsPullAndBear 9 "21:1"

$wsMango = $wb.Worksheets.Item("11 - 10 Mango")
Add-Row $wsMango 10 "21:1"
Add-Row $wsMango 11 "21:1"
Add-Row $wsMango 12 "21:2"
Add-Row $wsMango 13 "21:12"
Add-Row $wsMango 14 "21:35"
Add-Row $wsMango 15 "22:9"

$wsZara = $wb.Worksheets.Item("11 - 10 Zara")
Add-Row $wsZara 10 "21:1"
Add-Row $wsZara 11 "21:1"
Add-Row $wsZara 12 "21:2"
Add-Row $wsZara 13 "21:12"
Add-Row $wsZara 14 "21:35"
Add-Row $wsZara 15 "22:9"

$wsStradivarius = $wb.Worksheets.Item("11 - 10 Stradivarius")
Add-Row $wsStradivarius 10 "21:1"
Add-Row $wsStradivarius 11 "21:1"
Add-Row $wsStradivarius 12 "21:2"
Add-Row $wsStradivarius 13 "21:12"
Add-Row $wsStradivarius 14 "21:36"
Add-Row $wsStradivarius 15 "22:9"

$wsBershka = $wb.Worksheets.Item("11 - 10 Bershka")
Add-Row $wsBershka 10 "21:1"
Add-Row $wsBershka 11 "21:1"
Add-Row $wsBershka 12 "21:2"
Add-Row $wsBershka 13 "21:12"
Add-Row $wsBershka 14 "21:36"
Add-Row $wsBershka 15 "22:9"
